$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '43.972.10'
$ws.Range("E2").Value = '  +4.05%  '

$ws.Range("D3").Value = '2.209.74'
$ws.Range("E3").Value = '  +1.69%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '259.91'
$ws.Range("E5").Value = '  +2.68%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '81.93'

$ws.Range("E7").Value = '  +2.78%  '

$ws.Range("E8").Value = '  -0.15%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.602'
$ws.Range("E9").Value = '  +3.69%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.67'
$ws.Range("E10").Value = '  +7.50%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0927'
$ws.Range("E11").Value = '  +2.06%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.02'
$ws.Range("E12").Value = '  +4.12%  '

$ws.Range("E13").Value = '  +2.69%  '

$ws.Range("D14").Value = '2.541.39'
$ws.Range("E14").Value = '  +1.54%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.46'
$ws.Range("E15").Value = '  +2.46%  '

$ws.Range("D16").Value = '2.225.21'
$ws.Range("E16").Value = '  +1.92%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.779'
$ws.Range("E17").Value = '  +2.03%  '

$ws.Range("D18").Value = '43.842.34'
$ws.Range("E18").Value = '  +3.88%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000103'
$ws.Range("E19").Value = '  +1.41%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.13'
$ws.Range("E20").Value = '  +0.95%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.98'
$ws.Range("E21").Value = '  +2.27%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.38'
$ws.Range("E22").Value = '  +10.46%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.05'
$ws.Range("E23").Value = '  +2.64%  '

$ws.Range("E24").Value = '  -2.90%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  +0.06%  '

$ws.Range("E26").Value = '  +3.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '41.14'
$ws.Range("E27").Value = '  +10.77%  '

$ws.Range("E28").Value = '  +0.78%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.24'
$ws.Range("E29").Value = '  +2.68%  '

$ws.Range("E30").Value = '  +0.23%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '172.84'
$ws.Range("E31").Value = '  +2.30%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.55'
$ws.Range("E32").Value = '  +2.96%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0871'
$ws.Range("E33").Value = '  +8.19%  '

$ws.Range("E34").Value = '  +4.31%  '

$ws.Range("E35").Value = '  +7.56%  '

$ws.Range("E36").Value = '  +2.23%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.52'
$ws.Range("E37").Value = '  +7.32%  '

$ws.Range("E38").Value = '  +7.53%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '13.29'
$ws.Range("E39").Value = '  +11.99%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.94'
$ws.Range("E40").Value = '  +20.83%  '

$ws.Range("E41").Value = '  +3.15%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.53'
$ws.Range("E42").Value = '  +7.54%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '62.87'
$ws.Range("E43").Value = '  +6.15%  '

$ws.Range("E44").Value = '  +2.94%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '102.52'
$ws.Range("E45").Value = '  +0.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0987'
$ws.Range("E46").Value = '  +1.90%  '

$ws.Range("E47").Value = '  +0.73%  '

$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.12'
$ws.Range("E48").Value = '  +3.20%  '

$ws.Range("B49").Value = 'TrustWalletToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.17'
$ws.Range("E49").Value = '  +4.24%  '

$ws.Range("E50").Value = '  +28.65%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.442'
$ws.Range("E51").Value = '  -5.78%  '

